$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-08-28 14:45:01"
$zhcn.Range("H4").Value = "2016-08-28 14:44:55"
$zhcn.Range("K4").Value = "2016-08-28 14:45:27"
$dede.Range("K4").Value = "2016-08-28 14:45:34"
